$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 340, shifting existing rows 340-378 down to 341-379
$ws.Rows.Item(340).Insert()

# Populate the new row 340 with data (copy constant columns from the row above, and the
# new data values as per the updated dataset)
$ws.Cells.Item(340, 1).Value = 9
$ws.Cells.Item(340, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(340, 3).Value = "Metropolitana"
$ws.Cells.Item(340, 4).Value = 45124
$ws.Cells.Item(340, 5).Value = 13
$ws.Cells.Item(340, 6).Value = 100112001
$ws.Cells.Item(340, 7).Value = "Berenjena"
$ws.Cells.Item(340, 8).Value = "Sin especificar"
$ws.Cells.Item(340, 9).Value = "Primera"
$ws.Cells.Item(340, 10).Value = 70
$ws.Cells.Item(340, 11).Value = 6000
$ws.Cells.Item(340, 12).Value = 7000
$ws.Cells.Item(340, 13).Value = 6500
$ws.Cells.Item(340, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(340, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(340, 16).Value = 130
$ws.Cells.Item(340, 17).Value = 50
$ws.Cells.Item(340, 18).Value = "Hortaliza"
